# Changes in QUESTIONS table and data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions")

# Insert a new column before column C (shifts old C -> D, old D -> E)
$ws.Range("C1").EntireColumn.Insert()

# New column header
$ws.Range("C2").Value = "I_QSTN"

# Fill new column C (rows 3..116) with an incrementing counter 1..114
$n = 1
for ($r = 3; $r -le 116; $r++) {
    $ws.Cells.Item($r, 3).Value = $n
    $n = $n + 1
}

# Update the view: selected range / active cell and scroll position
$ws.Activate()
$ws.Range("A3:E65").Select()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 2
